$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Csf3"
$ws.Range("C2").Value = "Csf3r"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 1.414388666666667
$ws.Range("H2").Value = 4.243166
$ws.Range("I2").Value = 0.5586801917371232
$ws.Range("J2").Value = 0.5586801917371234
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 300.1573296666667
$ws.Range("N2").Value = 900.471989
$ws.Range("O2").Value = 0.8617605198060754
$ws.Range("P2").Value = 0.8617605198060753
$ws.Range("Q2").Value = 424.5391252974638
$ws.Range("R2").Value = 3820.852127677174
$ws.Range("S2").Value = 0.4814485324367412
$ws.Range("T2").Value = 0.4814485324367412

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Csf3"
$ws.Range("C3").Value = "Csf3r"
$ws.Range("D3").Value = "M2"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 1.414388666666667
$ws.Range("H3").Value = 4.243166
$ws.Range("I3").Value = 0.5586801917371232
$ws.Range("J3").Value = 0.5586801917371234
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 48.149796
$ws.Range("N3").Value = 144.449388
$ws.Range("O3").Value = 0.1382394801939247
$ws.Range("P3").Value = 0.1382394801939247
$ws.Range("Q3").Value = 68.10252576471201
$ws.Range("R3").Value = 612.922731882408
$ws.Range("S3").Value = 0.07723165930038209
$ws.Range("T3").Value = 0.07723165930038209

# Row 4
$ws.Range("A4").Value = "FAPs"
$ws.Range("B4").Value = "Csf3"
$ws.Range("C4").Value = "Csf3r"
$ws.Range("D4").Value = "ECs"
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 0.6471403333333333
$ws.Range("H4").Value = 1.941421
$ws.Range("I4").Value = 0.255618907325916
$ws.Range("J4").Value = 0.255618907325916
$ws.Range("K4").Value = 2
$ws.Range("L4").Value = 0.6666666666666666
$ws.Range("M4").Value = 300.1573296666667
$ws.Range("N4").Value = 900.471989
$ws.Range("O4").Value = 0.8617605198060754
$ws.Range("P4").Value = 0.8617605198060753
$ws.Range("Q4").Value = 194.2439143729299
$ws.Range("R4").Value = 1748.195229356369
$ws.Range("S4").Value = 0.2202822824494423
$ws.Range("T4").Value = 0.2202822824494423

# Row 5
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Csf3"
$ws.Range("C5").Value = "Csf3r"
$ws.Range("D5").Value = "M2"
$ws.Range("E5").Value = 2
$ws.Range("F5").Value = 0.6666666666666666
$ws.Range("G5").Value = 0.6471403333333333
$ws.Range("H5").Value = 1.941421
$ws.Range("I5").Value = 0.255618907325916
$ws.Range("J5").Value = 0.255618907325916
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 48.149796
$ws.Range("N5").Value = 144.449388
$ws.Range("O5").Value = 0.1382394801939247
$ws.Range("P5").Value = 0.1382394801939247
$ws.Range("Q5").Value = 31.159675033372
$ws.Range("R5").Value = 280.437075300348
$ws.Range("S5").Value = 0.03533662487647363
$ws.Range("T5").Value = 0.03533662487647362

# Row 6
$ws.Range("A6").Value = "M2"
$ws.Range("B6").Value = "Csf3"
$ws.Range("C6").Value = "Csf3r"
$ws.Range("D6").Value = "ECs"
$ws.Range("E6").Value = 2
$ws.Range("F6").Value = 0.6666666666666666
$ws.Range("G6").Value = 0.227516
$ws.Range("H6").Value = 0.6825479999999999
$ws.Range("I6").Value = 0.08986828408546589
$ws.Range("J6").Value = 0.08986828408546589
$ws.Range("K6").Value = 2
$ws.Range("L6").Value = 0.6666666666666666
$ws.Range("M6").Value = 300.1573296666667
$ws.Range("N6").Value = 900.471989
$ws.Range("O6").Value = 0.8617605198060754
$ws.Range("P6").Value = 0.8617605198060753
$ws.Range("Q6").Value = 68.29059501644133
$ws.Range("R6").Value = 614.615355147972
$ws.Range("S6").Value = 0.07744493920757113
$ws.Range("T6").Value = 0.07744493920757113

# Row 7
$ws.Range("A7").Value = "M2"
$ws.Range("B7").Value = "Csf3"
$ws.Range("C7").Value = "Csf3r"
$ws.Range("D7").Value = "M2"
$ws.Range("E7").Value = 2
$ws.Range("F7").Value = 0.6666666666666666
$ws.Range("G7").Value = 0.227516
$ws.Range("H7").Value = 0.6825479999999999
$ws.Range("I7").Value = 0.08986828408546589
$ws.Range("J7").Value = 0.08986828408546589
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 48.149796
$ws.Range("N7").Value = 144.449388
$ws.Range("O7").Value = 0.1382394801939247
$ws.Range("P7").Value = 0.1382394801939247
$ws.Range("Q7").Value = 10.954848986736
$ws.Range("R7").Value = 98.593640880624
$ws.Range("S7").Value = 0.01242334487789476
$ws.Range("T7").Value = 0.01242334487789476

# Row 8
$ws.Range("A8").Value = "sCs"
$ws.Range("B8").Value = "Csf3"
$ws.Range("C8").Value = "Csf3r"
$ws.Range("D8").Value = "ECs"
$ws.Range("E8").Value = 1
$ws.Range("F8").Value = 0.3333333333333333
$ws.Range("G8").Value = 0.2426156666666667
$ws.Range("H8").Value = 0.727847
$ws.Range("I8").Value = 0.09583261685149484
$ws.Range("J8").Value = 0.09583261685149484
$ws.Range("K8").Value = 2
$ws.Range("L8").Value = 0.6666666666666666
$ws.Range("M8").Value = 300.1573296666667
$ws.Range("N8").Value = 900.471989
$ws.Range("O8").Value = 0.8617605198060754
$ws.Range("P8").Value = 0.8617605198060753
$ws.Range("Q8").Value = 72.82287064196478
$ws.Range("R8").Value = 655.405835777683
$ws.Range("S8").Value = 0.08258476571232065
$ws.Range("T8").Value = 0.08258476571232064

# Row 9
$ws.Range("A9").Value = "sCs"
$ws.Range("B9").Value = "Csf3"
$ws.Range("C9").Value = "Csf3r"
$ws.Range("D9").Value = "M2"
$ws.Range("E9").Value = 1
$ws.Range("F9").Value = 0.3333333333333333
$ws.Range("G9").Value = 0.2426156666666667
$ws.Range("H9").Value = 0.727847
$ws.Range("I9").Value = 0.09583261685149484
$ws.Range("J9").Value = 0.09583261685149484
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 48.149796
$ws.Range("N9").Value = 144.449388
$ws.Range("O9").Value = 0.1382394801939247
$ws.Range("P9").Value = 0.1382394801939247
$ws.Range("Q9").Value = 11.681894856404
$ws.Range("R9").Value = 105.137053707636
$ws.Range("S9").Value = 0.01324785113917419
$ws.Range("T9").Value = 0.01324785113917419
